$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated vm_pu values for the "case with 380 kV" run (rows 2-25, columns B-F and I-N)
$data = @{
  2 = @{ "B"="1.02"; "C"="1.035270263719528"; "D"="1.043076811231823"; "E"="1.044183789309923"; "F"="1.054714861396286"; "I"="1.038397002444473"; "J"="1.040384955810975"; "K"="1.045851695884775"; "L"="1.04695555494004"; "M"="1.057457309517885"; "N"="1.017423131858447" }
  3 = @{ "B"="1.02"; "C"="1.036026823282825"; "D"="1.0436582246267"; "E"="1.04485376362757"; "F"="1.055458533832535"; "I"="1.038541306548671"; "J"="1.04078596820538"; "K"="1.046244671670251"; "L"="1.047437084379395"; "M"="1.058014457932326"; "N"="1.017556977977764" }
  4 = @{ "B"="1.02"; "C"="1.036517165471297"; "D"="1.044035148595964"; "E"="1.045288391285698"; "F"="1.055940977179822"; "I"="1.038633874519091"; "J"="1.041045546737753"; "K"="1.046498961626372"; "L"="1.047749091325752"; "M"="1.058375537567017"; "N"="1.017643591754525" }
  5 = @{ "B"="1.02"; "C"="1.036723494237189"; "D"="1.044193775550109"; "E"="1.045471372245352"; "F"="1.056144090355706"; "I"="1.038672596258484"; "J"="1.041154695314578"; "K"="1.046605865677588"; "L"="1.047880359225302"; "M"="1.05852746942932"; "N"="1.017680005219216" }
  6 = @{ "B"="1.02"; "C"="1.036758148763701"; "D"="1.044220419513746"; "E"="1.045502110943282"; "F"="1.0561782111252"; "I"="1.038679086419855"; "J"="1.041173023076654"; "K"="1.046623815328639"; "L"="1.047902405510827"; "M"="1.058552987260812"; "N"="1.017686119246295" }
  7 = @{ "B"="1.02"; "C"="1.03651992170662"; "D"="1.044037267519032"; "E"="1.045290835253561"; "F"="1.0559436900357"; "I"="1.038634392683472"; "J"="1.041047005102613"; "K"="1.046500390083176"; "L"="1.047750844942245"; "M"="1.058377567163849"; "N"="1.017644078309588" }
  8 = @{ "B"="1.02"; "C"="1.035525780206116"; "D"="1.04327315426521"; "E"="1.044409979270589"; "F"="1.054965931364327"; "I"="1.038445937204321"; "J"="1.040520458970934"; "K"="1.045984501295392"; "L"="1.047118201057879"; "M"="1.057645481851472"; "N"="1.017468364160367" }
  9 = @{ "B"="1.02"; "C"="1.033780174041983"; "D"="1.041932220363731"; "E"="1.042866389400014"; "F"="1.053252574684814"; "I"="1.038107712527219"; "J"="1.039593419408375"; "K"="1.045075567141974"; "L"="1.046006729633018"; "M"="1.056359880735863"; "N"="1.01715880424119" }
  10 = @{ "B"="1.02"; "C"="1.032620725441579"; "D"="1.041042103555988"; "E"="1.041843232416281"; "F"="1.052116913913111"; "I"="1.037878146288856"; "J"="1.03897601538743"; "K"="1.044469783490319"; "L"="1.045268082899206"; "M"="1.055505899978487"; "N"="1.016952508307263" }
  11 = @{ "B"="1.02"; "C"="1.032119714435861"; "D"="1.040657609288202"; "E"="1.041401621215813"; "F"="1.051626749284946"; "I"="1.037777782966066"; "J"="1.038708837893024"; "K"="1.044207531757733"; "L"="1.044948814895993"; "M"="1.055136872445136"; "N"="1.016863204539522" }
  12 = @{ "B"="1.02"; "C"="1.031933774670173"; "D"="1.040514933137131"; "E"="1.041237803081656"; "F"="1.051444920672018"; "I"="1.037740360175867"; "J"="1.038609622127537"; "K"="1.044110129782944"; "L"="1.044830311992619"; "M"="1.054999914048072"; "N"="1.01683003722202" }
  13 = @{ "B"="1.02"; "C"="1.031973652196241"; "D"="1.040545531190243"; "E"="1.041272932820868"; "F"="1.051483912612006"; "I"="1.037748393967113"; "J"="1.038630903053144"; "K"="1.044131022360889"; "L"="1.044855727283548"; "M"="1.055029286860794"; "N"="1.016837151531451" }
  14 = @{ "B"="1.02"; "C"="1.032104341360936"; "D"="1.040645812721375"; "E"="1.041388075540099"; "F"="1.051611714358172"; "I"="1.037774692506962"; "J"="1.03870063615123"; "K"="1.044199480267573"; "L"="1.044939017615427"; "M"="1.055125549067142"; "N"="1.016860462832558" }
  15 = @{ "B"="1.02"; "C"="1.032184884271791"; "D"="1.040707618379709"; "E"="1.041459047448253"; "F"="1.051690489152393"; "I"="1.037790876929268"; "J"="1.03874360448299"; "K"="1.044241660812432"; "L"="1.044990347175685"; "M"="1.05518487466679"; "N"="1.016874826250241" }
  16 = @{ "B"="1.02"; "C"="1.032653997938204"; "D"="1.041067640992744"; "E"="1.041872570856966"; "F"="1.052149478100998"; "I"="1.037884786921193"; "J"="1.038993750624462"; "K"="1.044487189593876"; "L"="1.045289283836573"; "M"="1.055530407122951"; "N"="1.016958435648002" }
  17 = @{ "B"="1.02"; "C"="1.032948539900634"; "D"="1.041293724714732"; "E"="1.042132345577972"; "F"="1.052437815499435"; "I"="1.03794343781577"; "J"="1.03915070521515"; "K"="1.044641219516568"; "L"="1.045476952938062"; "M"="1.055747353126555"; "N"="1.017010888330298" }
  18 = @{ "B"="1.02"; "C"="1.033120441244576"; "D"="1.041425685327061"; "E"="1.04228400488276"; "F"="1.052606150497396"; "I"="1.037977555247201"; "J"="1.039242269757171"; "K"="1.044731067924124"; "L"="1.045586472121796"; "M"="1.055873966466643"; "N"="1.017041485348015" }
  19 = @{ "B"="1.02"; "C"="1.033179072058535"; "D"="1.041470695654098"; "E"="1.042335739974567"; "F"="1.052663574204869"; "I"="1.037989172671768"; "J"="1.039273493499525"; "K"="1.044761704772821"; "L"="1.045623824606392"; "M"="1.05591715056676"; "N"="1.017051918507339" }
  20 = @{ "B"="1.02"; "C"="1.032916927976941"; "D"="1.041269458767257"; "E"="1.042104460011539"; "F"="1.052406863816675"; "I"="1.037937154706671"; "J"="1.039133863850626"; "K"="1.044624692998448"; "L"="1.045456812121466"; "M"="1.055724069370365"; "N"="1.017005260416896" }
  21 = @{ "B"="1.02"; "C"="1.032065852280038"; "D"="1.040616278369403"; "E"="1.041354162900734"; "F"="1.05157407327525"; "I"="1.037766952194188"; "J"="1.038680100766737"; "K"="1.044179320832411"; "L"="1.044914488261193"; "M"="1.05509719906059"; "N"="1.016853598118715" }
  22 = @{ "B"="1.02"; "C"="1.03153166255873"; "D"="1.040206421335486"; "E"="1.040883671714228"; "F"="1.051051856600938"; "I"="1.037659109938643"; "J"="1.038394952187277"; "K"="1.043899356243805"; "L"="1.044574014503059"; "M"="1.054703726305567"; "N"="1.016758265860034" }
  23 = @{ "B"="1.02"; "C"="1.03181475917185"; "D"="1.040423615557918"; "E"="1.041132968735131"; "F"="1.051328560821249"; "I"="1.037716357491509"; "J"="1.03854610012066"; "K"="1.044047764719217"; "L"="1.044754457488495"; "M"="1.054912249938369"; "N"="1.016808800862712" }
  24 = @{ "B"="1.02"; "C"="1.032931211735838"; "D"="1.04128042322722"; "E"="1.042117059875365"; "F"="1.052420849078414"; "I"="1.037939994059579"; "J"="1.039141473690368"; "K"="1.044632160604255"; "L"="1.04546591271991"; "M"="1.055734590074089"; "N"="1.017007803421813" }
  25 = @{ "B"="1.02"; "C"="1.034230708043052"; "D"="1.042278216033847"; "E"="1.043264413226328"; "F"="1.053694369203349"; "I"="1.03819587479693"; "J"="1.039832977745634"; "K"="1.045310524342418"; "L"="1.046293667289914"; "M"="1.056691703468877"; "N"="1.017238821439941" }
}

foreach ($r in $data.Keys) {
  $rowData = $data[$r]
  foreach ($c in $rowData.Keys) {
    $ws.Range("$c$r").Value = [double]$rowData[$c]
  }
}
